$d = $word.ActiveDocument

# 1. Title meta content: "1 Introduction to Game Maker" -> "8 The Maze"
$d.Content.Find.Execute("1 Introduction to Game Maker", $false, $false, $false, $false, $false,
                         $true, 1, $false, "8 The Maze", 2) | Out-Null

# 2. Description meta content: collapse highlighted "This " + "is our first article..." + ". "
#    into a single plain (non-highlighted) run with new text. Find/Replace would inherit the
#    formatting (yellow highlight) of the first matched run, but the target has no rPr at all,
#    so splice the paragraph's OOXML directly (preserving the surrounding runs/rsids untouched).
$descPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*is our first article*") {
        $descPara = $para
        break
    }
}
if ($descPara -ne $null) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="135F4BCC" w14:textId="1E58D446" w:rsidR="00B876FE" w:rsidRPr="009869DD" w:rsidRDefault="00B876FE" w:rsidP="00B876FE"><w:r w:rsidRPr="009869DD"><w:t xml:space="preserve">&lt;meta name="</w:t></w:r><w:r w:rsidRPr="009869DD"><w:rPr><w:b/><w:i/><w:color w:val="3333FF"/></w:rPr><w:t>description</w:t></w:r><w:r w:rsidRPr="009869DD"><w:t xml:space="preserve">" content="</w:t></w:r><w:r><w:t>In this article, we will be looking at how to create a maze for our room.</w:t></w:r><w:r w:rsidRPr="009869DD"><w:t>"/&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $descPara.Range.InsertXML($xml)
}

# 3. Remove the grammar-check proofing marks (gramStart/gramEnd) around "index,follow",
#    while leaving the spelling marks (spellStart/spellEnd) untouched. The Word object model
#    has no direct "remove this proofing mark" call, so splice the paragraph's OOXML back in
#    without the gramStart/gramEnd elements (everything else - text, rPr, rsids - is preserved).
$robotsPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*index,follow*") {
        $robotsPara = $para
        break
    }
}
if ($robotsPara -ne $null) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4A2F2627" w14:textId="77777777" w:rsidR="00B876FE" w:rsidRPr="009869DD" w:rsidRDefault="00B876FE" w:rsidP="00B876FE"><w:r w:rsidRPr="009869DD"><w:t xml:space="preserve">&lt;meta name="</w:t></w:r><w:r w:rsidRPr="009869DD"><w:rPr><w:b/><w:i/><w:color w:val="3333FF"/></w:rPr><w:t>robots</w:t></w:r><w:r w:rsidRPr="009869DD"><w:t xml:space="preserve">" content="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009869DD"><w:t>index,follow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009869DD"><w:t xml:space="preserve">" /&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $robotsPara.Range.InsertXML($xml)
}

# 4. Revised date meta content: "Thursday 30" -> "November, 28"
$d.Content.Find.Execute("Thursday 30", $false, $false, $false, $false, $false,
                         $true, 1, $false, "November, 28", 2) | Out-Null

# 5. URL meta content: swap in the new article path.
$d.Content.Find.Execute("I-Snuck-A-Book/PDF_Optimizer.html", $false, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Enlightenment/Articles/2025/4_Game_Maker/8_The_Maze/8_The_Maze.html", 2) | Out-Null
